# Pitch opening slide: change the highlighted "DISCIPLINA" title block
# from yellow (FFFF00) to red (FF0000) highlight color.
#
# Commit message: "Alteracao de fundo da pagina de abertura do Pitch,
# com a cor vermelha" (change of the opening-page highlight to red).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Find the "CaixaDeTexto 16" textbox that holds the DISCIPLINA / PROJETO
# DE SISTEMAS ... / QUALIDADE DE SOFTWARE E GOVERNANCA DE TI text, which
# is the only shape using the yellow highlight on this slide.
$shp = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Name -eq "CaixaDeTexto 16") {
        $shp = $candidate
    }
}
if ($shp -eq $null) {
    # Fallback: it is the last shape on the slide.
    $shp = $s.Shapes.Item($s.Shapes.Count)
}

$tr = $shp.TextFrame.TextRange

# New highlight color: red (RGB 255,0,0 -> COM RGB() = 255).
$newHighlight = 255

# Run 1 (paragraph 1): "DISCIPLINA:   "
$tr.Characters(1, 14).Font.Highlight.RGB = $newHighlight

# Run 2 (paragraph 1): "PROJETO DE SISTEMAS APLICADO AS MELHORES PRATICAS EM "
$tr.Characters(15, 53).Font.Highlight.RGB = $newHighlight

# Run 3 (paragraph 2): "QUALIDADE DE SOFTWARE E GOVERNANCA DE TI"
$tr.Characters(69, 40).Font.Highlight.RGB = $newHighlight

# Paragraph-end mark of paragraph 2 (maps to <a:endParaRPr> in the XML)
# also carries the same yellow highlight in the source file; update it too
# for completeness.
$tr.Paragraphs(2, 1).Font.Highlight.RGB = $newHighlight
